$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.08975000000000001
$ws.Range("E2").Value = -0.08259999999999999
$ws.Range("G2").Value = 0.0427957605270696
$ws.Range("H2").Value = 0.0427957605270696
$ws.Range("I2").Value = 0.04219421369235176
$ws.Range("J2").Value = 0.03265898587368644
$ws.Range("K2").Value = 4.390000000000001
$ws.Range("L2").Value = 0.01257519335433973
$ws.Range("M2").Value = 6.37
$ws.Range("N2").Value = 0.04125647668393782
$ws.Range("O2").Value = 1.451025056947608
$ws.Range("P2").Value = 4.67
$ws.Range("Q2").Value = 0.0302461139896373
$ws.Range("R2").Value = 1.06378132118451
$ws.Range("S2").Value = 1.7
$ws.Range("T2").Value = 0.2668759811616954
$ws.Range("U2").Value = 78.40000000000001
$ws.Range("V2").Value = 0.5077720207253886
$ws.Range("W2").Value = 0.006792812059882114
$ws.Range("X2").Value = 0.0669874317059061
$ws.Range("Y2").Value = -0.06019461964602398
$ws.Range("Z2").Value = 1.652724320281024
$ws.Range("AA2").Value = 0.04101944812220443
$ws.Range("AB2").Value = 0.06612953444188396
$ws.Range("AC2").Value = -0.02511008631967953
$ws.Range("AD2").Value = 2.536
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2.536
$ws.Range("AG2").Value = -75.864
$ws.Range("AH2").Value = 0.01615945353519906
$ws.Range("AI2").Value = 0.007528886461067107
$ws.Range("AJ2").Value = -0.9659773861668535
$ws.Range("AK2").Value = -0.2935504341500411
$ws.Range("AL2").Value = 0.151
$ws.Range("AM2").Value = 0.151
$ws.Range("AN2").Value = 0.1436827195467422
$ws.Range("AO2").Value = 97.5496688741722
$ws.Range("AP2").Value = -4.298243626062323
$ws.Range("AQ2").Value = 97.5496688741722
$ws.Range("D3").Value = 0.09390000000000001
$ws.Range("E3").Value = -0.08259999999999999
$ws.Range("G3").Value = 0.05621970920840064
$ws.Range("H3").Value = 0.05621970920840064
$ws.Range("I3").Value = 0.04135702746365105
$ws.Range("J3").Value = 0.02266495363362294
$ws.Range("K3").Value = 6.99
$ws.Range("L3").Value = 0.02258481421647819
$ws.Range("M3").Value = 4.67
$ws.Range("N3").Value = 0.0451207729468599
$ws.Range("O3").Value = 0.6680972818311873
$ws.Range("P3").Value = 4.67
$ws.Range("Q3").Value = 0.0451207729468599
$ws.Range("R3").Value = 0.6680972818311873
$ws.Range("U3").Value = 44.3
$ws.Range("V3").Value = 0.4280193236714975
$ws.Range("W3").Value = 0.0346382556987116
$ws.Range("X3").Value = 0.06629608395392414
$ws.Range("Y3").Value = -0.03165782825521254
$ws.Range("Z3").Value = 1.921322018536567
$ws.Range("AA3").Value = 0.04354667446539012
$ws.Range("AB3").Value = 0.06604202419794314
$ws.Range("AC3").Value = -0.02249534973255302
$ws.Range("AD3").Value = 0.606
$ws.Range("AF3").Value = 0.606
$ws.Range("AG3").Value = -43.694
$ws.Range("AH3").Value = 0.005820990144660251
$ws.Range("AI3").Value = 0.002854370578316204
$ws.Range("AJ3").Value = -0.7305955924154766
$ws.Range("AK3").Value = -0.2600740449745842
$ws.Range("AN3").Value = 0.04067114093959731
$ws.Range("AO3").Value = 3200
$ws.Range("AP3").Value = -2.93248322147651
$ws.Range("AQ3").Value = 3200
$ws.Range("D4").Value = 0.08560000000000001
$ws.Range("G4").Value = -0.06212121212121212
$ws.Range("H4").Value = -0.06212121212121212
$ws.Range("I4").Value = 0.04873737373737373
$ws.Range("J4").Value = 0.04873737373737373
$ws.Range("K4").Value = -2.6
$ws.Range("L4").Value = -0.06565656565656566
$ws.Range("M4").Value = 1.7
$ws.Range("N4").Value = 0.03339882121807466
$ws.Range("O4").Value = -0.6538461538461539
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 1.7
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 34.1
$ws.Range("V4").Value = 0.6699410609037328
$ws.Range("W4").Value = -0.02105263157894737
$ws.Range("X4").Value = 0.06767877945788806
$ws.Range("Y4").Value = -0.08873141103683542
$ws.Range("Z4").Value = 0.7897885919425609
$ws.Range("AA4").Value = 0.03849222177901875
$ws.Range("AB4").Value = 0.06621704468582477
$ws.Range("AC4").Value = -0.02772482290680603
$ws.Range("AD4").Value = 1.93
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.93
$ws.Range("AG4").Value = -32.17
$ws.Range("AH4").Value = 0.03653227332954761
$ws.Range("AI4").Value = 0.01549827350839155
$ws.Range("AJ4").Value = -1.717565403096637
$ws.Range("AK4").Value = -0.3557447749640607
$ws.Range("AL4").Value = 0.147
$ws.Range("AM4").Value = 0.147
$ws.Range("AN4").Value = 0.7018181818181818
$ws.Range("AO4").Value = 13.12925170068027
$ws.Range("AP4").Value = -11.69818181818182
$ws.Range("AQ4").Value = 13.12925170068027
